$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Fgf1"
$ws.Cells.Item(2,3).Value = "Fgfr3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.190640666666667
$ws.Cells.Item(2,8).Value = 3.571922
$ws.Cells.Item(2,9).Value = 0.1136540143525372
$ws.Cells.Item(2,10).Value = 0.1136540143525372
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 5.751166666666666
$ws.Cells.Item(2,14).Value = 17.2535
$ws.Cells.Item(2,15).Value = 0.7405222614421495
$ws.Cells.Item(2,16).Value = 0.7405222614421495
$ws.Cells.Item(2,17).Value = 6.847572914111111
$ws.Cells.Item(2,18).Value = 61.628156227
$ws.Cells.Item(2,19).Value = 0.0841633277303194
$ws.Cells.Item(2,20).Value = 0.0841633277303194

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Fgf1"
$ws.Cells.Item(3,3).Value = "Fgfr3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.190640666666667
$ws.Cells.Item(3,8).Value = 3.571922
$ws.Cells.Item(3,9).Value = 0.1136540143525372
$ws.Cells.Item(3,10).Value = 0.1136540143525372
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.5698483333333333
$ws.Cells.Item(3,14).Value = 1.709545
$ws.Cells.Item(3,15).Value = 0.07337387367415998
$ws.Cells.Item(3,16).Value = 0.07337387367416
$ws.Cells.Item(3,17).Value = 0.6784845994988888
$ws.Cells.Item(3,18).Value = 6.10636139549
$ws.Cells.Item(3,19).Value = 0.008339235291664233
$ws.Cells.Item(3,20).Value = 0.008339235291664235

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Fgf1"
$ws.Cells.Item(4,3).Value = "Fgfr3"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.190640666666667
$ws.Cells.Item(4,8).Value = 3.571922
$ws.Cells.Item(4,9).Value = 0.1136540143525372
$ws.Cells.Item(4,10).Value = 0.1136540143525372
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.445350666666667
$ws.Cells.Item(4,14).Value = 4.336052
$ws.Cells.Item(4,15).Value = 0.1861038648836906
$ws.Cells.Item(4,16).Value = 0.1861038648836906
$ws.Cells.Item(4,17).Value = 1.720893281327111
$ws.Cells.Item(4,18).Value = 15.488039531944
$ws.Cells.Item(4,19).Value = 0.02115145133055362
$ws.Cells.Item(4,20).Value = 0.02115145133055362

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Fgf1"
$ws.Cells.Item(5,3).Value = "Fgfr3"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 5.168173666666667
$ws.Cells.Item(5,8).Value = 15.504521
$ws.Cells.Item(5,9).Value = 0.4933341355895272
$ws.Cells.Item(5,10).Value = 0.4933341355895272
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 5.751166666666666
$ws.Cells.Item(5,14).Value = 17.2535
$ws.Cells.Item(5,15).Value = 0.7405222614421495
$ws.Cells.Item(5,16).Value = 0.7405222614421495
$ws.Cells.Item(5,17).Value = 29.72302811927778
$ws.Cells.Item(5,18).Value = 267.5072530735
$ws.Cells.Item(5,19).Value = 0.3653249097333647
$ws.Cells.Item(5,20).Value = 0.3653249097333647

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Fgf1"
$ws.Cells.Item(6,3).Value = "Fgfr3"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 5.168173666666667
$ws.Cells.Item(6,8).Value = 15.504521
$ws.Cells.Item(6,9).Value = 0.4933341355895272
$ws.Cells.Item(6,10).Value = 0.4933341355895272
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.5698483333333333
$ws.Cells.Item(6,14).Value = 1.709545
$ws.Cells.Item(6,15).Value = 0.07337387367415998
$ws.Cells.Item(6,16).Value = 0.07337387367416
$ws.Cells.Item(6,17).Value = 2.945075150327222
$ws.Cells.Item(6,18).Value = 26.505676352945
$ws.Cells.Item(6,19).Value = 0.03619783654389688
$ws.Cells.Item(6,20).Value = 0.03619783654389689

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Fgf1"
$ws.Cells.Item(7,3).Value = "Fgfr3"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 5.168173666666667
$ws.Cells.Item(7,8).Value = 15.504521
$ws.Cells.Item(7,9).Value = 0.4933341355895272
$ws.Cells.Item(7,10).Value = 0.4933341355895272
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.445350666666667
$ws.Cells.Item(7,14).Value = 4.336052
$ws.Cells.Item(7,15).Value = 0.1861038648836906
$ws.Cells.Item(7,16).Value = 0.1861038648836906
$ws.Cells.Item(7,17).Value = 7.46982325456578
$ws.Cells.Item(7,18).Value = 67.228409291092
$ws.Cells.Item(7,19).Value = 0.09181138931226566
$ws.Cells.Item(7,20).Value = 0.09181138931226566

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Fgf1"
$ws.Cells.Item(8,3).Value = "Fgfr3"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 4.037194
$ws.Cells.Item(8,8).Value = 12.111582
$ws.Cells.Item(8,9).Value = 0.3853751326204581
$ws.Cells.Item(8,10).Value = 0.3853751326204581
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 5.751166666666666
$ws.Cells.Item(8,14).Value = 17.2535
$ws.Cells.Item(8,15).Value = 0.7405222614421495
$ws.Cells.Item(8,16).Value = 0.7405222614421495
$ws.Cells.Item(8,17).Value = 23.21857555966667
$ws.Cells.Item(8,18).Value = 208.967180037
$ws.Cells.Item(8,19).Value = 0.2853788647116699
$ws.Cells.Item(8,20).Value = 0.2853788647116699

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Fgf1"
$ws.Cells.Item(9,3).Value = "Fgfr3"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 4.037194
$ws.Cells.Item(9,8).Value = 12.111582
$ws.Cells.Item(9,9).Value = 0.3853751326204581
$ws.Cells.Item(9,10).Value = 0.3853751326204581
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.5698483333333333
$ws.Cells.Item(9,14).Value = 1.709545
$ws.Cells.Item(9,15).Value = 0.07337387367415998
$ws.Cells.Item(9,16).Value = 0.07337387367416
$ws.Cells.Item(9,17).Value = 2.300588272243334
$ws.Cells.Item(9,18).Value = 20.70529445019
$ws.Cells.Item(9,19).Value = 0.02827646629805614
$ws.Cells.Item(9,20).Value = 0.02827646629805615

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Fgf1"
$ws.Cells.Item(10,3).Value = "Fgfr3"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 4.037194
$ws.Cells.Item(10,8).Value = 12.111582
$ws.Cells.Item(10,9).Value = 0.3853751326204581
$ws.Cells.Item(10,10).Value = 0.3853751326204581
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.445350666666667
$ws.Cells.Item(10,14).Value = 4.336052
$ws.Cells.Item(10,15).Value = 0.1861038648836906
$ws.Cells.Item(10,16).Value = 0.1861038648836906
$ws.Cells.Item(10,17).Value = 5.835161039362668
$ws.Cells.Item(10,18).Value = 52.516449354264
$ws.Cells.Item(10,19).Value = 0.07171980161073208
$ws.Cells.Item(10,20).Value = 0.07171980161073208

# Row 11
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,2).Value = "Fgf1"
$ws.Cells.Item(11,3).Value = "Fgfr3"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 1
$ws.Cells.Item(11,6).Value = 0.3333333333333333
$ws.Cells.Item(11,7).Value = 0.08000233333333333
$ws.Cells.Item(11,8).Value = 0.240007
$ws.Cells.Item(11,9).Value = 0.007636717437477471
$ws.Cells.Item(11,10).Value = 0.007636717437477472
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 5.751166666666666
$ws.Cells.Item(11,14).Value = 17.2535
$ws.Cells.Item(11,15).Value = 0.7405222614421495
$ws.Cells.Item(11,16).Value = 0.7405222614421495
$ws.Cells.Item(11,17).Value = 0.4601067527222222
$ws.Cells.Item(11,18).Value = 4.1409607745
$ws.Cells.Item(11,19).Value = 0.005655159266795514
$ws.Cells.Item(11,20).Value = 0.005655159266795514

# Row 12
$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,2).Value = "Fgf1"
$ws.Cells.Item(12,3).Value = "Fgfr3"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 1
$ws.Cells.Item(12,6).Value = 0.3333333333333333
$ws.Cells.Item(12,7).Value = 0.08000233333333333
$ws.Cells.Item(12,8).Value = 0.240007
$ws.Cells.Item(12,9).Value = 0.007636717437477471
$ws.Cells.Item(12,10).Value = 0.007636717437477472
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 0.5698483333333333
$ws.Cells.Item(12,14).Value = 1.709545
$ws.Cells.Item(12,15).Value = 0.07337387367415998
$ws.Cells.Item(12,16).Value = 0.07337387367416
$ws.Cells.Item(12,17).Value = 0.04558919631277777
$ws.Cells.Item(12,18).Value = 0.4103027668149999
$ws.Cells.Item(12,19).Value = 0.0005603355405427267
$ws.Cells.Item(12,20).Value = 0.0005603355405427269

# Row 13
$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,2).Value = "Fgf1"
$ws.Cells.Item(13,3).Value = "Fgfr3"
$ws.Cells.Item(13,4).Value = "MuSCs"
$ws.Cells.Item(13,5).Value = 1
$ws.Cells.Item(13,6).Value = 0.3333333333333333
$ws.Cells.Item(13,7).Value = 0.08000233333333333
$ws.Cells.Item(13,8).Value = 0.240007
$ws.Cells.Item(13,9).Value = 0.007636717437477471
$ws.Cells.Item(13,10).Value = 0.007636717437477472
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 1.445350666666667
$ws.Cells.Item(13,14).Value = 4.336052
$ws.Cells.Item(13,15).Value = 0.1861038648836906
$ws.Cells.Item(13,16).Value = 0.1861038648836906
$ws.Cells.Item(13,17).Value = 0.1156314258182222
$ws.Cells.Item(13,18).Value = 1.040682832364
$ws.Cells.Item(13,19).Value = 0.001421222630139231
$ws.Cells.Item(13,20).Value = 0.001421222630139231
